# Update the "education" worksheet with the new degrees (PhD/MSc -> PhD/MSc in
# Neurociencias, Universidad de Valencia / Universidad Católica de Colombia),
# matching the author's latest CV data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# --- new string values (written in the same order they were first entered so
#     the shared-string table is rebuilt in the same sequence Excel used) ---
$valenciaLink  = "\href{https://www.uv.es/uvweb/universidad/es/universidad-valencia-1285845048380.html}{Universidad de Valencia}"
$valencia      = "Valencia, España"
$catolicaLink  = "\href{https://www.ucatolica.edu.co/portal/Pregrado/psicologia/}{Universidad Cátolica de Colombia}"
$bogota        = "Bogotá, Colombia"
$phdNeuro      = "PhD - Neurociencias"
$masterNeuro   = "Máster en Neurociencias Básicas y Aplicadas "
$psicologia    = "Psicología"

$ws.Range("C2").Value = $valenciaLink
$ws.Range("D2").Value = $valencia
$ws.Range("C3").Value = $valenciaLink
$ws.Range("D3").Value = $valencia
$ws.Range("C4").Value = $catolicaLink
$ws.Range("D4").Value = $bogota
$ws.Range("A2").Value = $phdNeuro
$ws.Range("A3").Value = $masterNeuro
$ws.Range("A4").Value = $psicologia

$ws.Range("B2").Value = 2018
$ws.Range("B3").Value = 2012
$ws.Range("B4").Value = 2007

# --- formatting: rows 2-4 (A:D) now wrap, left/top aligned, and are taller ---
$dataRange = $ws.Range("A2:D4")
$dataRange.WrapText = $true
$dataRange.HorizontalAlignment = -4131   # xlLeft
$dataRange.VerticalAlignment = -4160     # xlTop

$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8

# --- cursor left where the author last clicked before saving ---
$ws.Range("A22").Select()
